$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.947.12"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.817.92"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'310.02"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.4651"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("D9").Value = "'0.07357"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").Value = "'0.8717"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").Value = "'20.24"
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "1.820.80"
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("D13").Value = "'5.401"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").Value = "'0.07106"
$ws.Range("E14").Value = "  +0.86%  "
$ws.Range("D15").Value = "'6.505"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "'91.47"
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "'0.000008708"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("D21").Value = "26.958.44"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("D22").Value = "'5.294"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").Value = "'10.59"
$ws.Range("E23").Value = "  +0.12%  "
$ws.Range("D24").Value = "2.025.30"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "'150.63"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").Value = "'18.35"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").Value = "'5.246"
$ws.Range("E29").Value = "  -1.56%  "
$ws.Range("D30").Value = "'116.53"
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").Value = "'0.08898"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").Value = "'0.7577"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("D34").Value = "'4.500"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").Value = "'2.908"
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").Value = "'1.090"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").Value = "'0.05288"
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("D40").Value = "'2.974"
$ws.Range("E40").Value = "  +1.60%  "
$ws.Range("D41").Value = "'7.175"
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'0.5281"
$ws.Range("E42").Value = "  -0.76%  "
$ws.Range("D43").Value = "'2.347"
$ws.Range("E43").Value = "  -3.05%  "
$ws.Range("D44").Value = "'0.1656"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").Value = "'8.443"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").Value = "'0.4873"
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "'103.29"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").Value = "'0.06293"
$ws.Range("E51").Value = "  +0.00%  "
